$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact / No display for ContactDetail" row (row 11)
$meta.Rows.Item(11).Delete()

# Row 10 (previously the first "Contact" row) becomes "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root extension row: Short / Definition now reflect the real title & description
$elements.Range("K2").Value = "Coverage Insurance Plan"
$elements.Range("L2").Value = "Reference to the insurance plan for this coverage"
